$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 view tweaks -------------------------------------------------
# Select C32 on Sheet1 while it is still the active sheet, so its
# sheetView ends up with a plain <selection> (no tabSelected) once Sheet2
# becomes the active tab below.
$ws1.Range("C32").Select()

# Column A gets wider.
$ws1.Columns.Item(1).ColumnWidth = 20.15

# --- Add Sheet2 right after Sheet1 --------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Column A on the new sheet.
$ws2.Columns.Item(1).ColumnWidth = 28.7

# --- Header row ----------------------------------------------------------
$ws2.Range("A1").Value = "Name"
$ws2.Range("A1").Font.Name = "Arial"
$ws2.Range("A1").Interior.ColorIndex = 26
$ws2.Range("A1").HorizontalAlignment = -4108

$ws2.Range("B1").Value = "Hmhw"
$ws2.Range("C1").Value = "Hmtl"
$ws2.Range("D1").Value = "Hmlw"
$ws2.Range("B1:D1").HorizontalAlignment = -4108

# --- Data rows -------------------------------------------------------------
$names = @(
    "Hayle Estuary",
    "Gannel Estuary",
    "Camel Estuary",
    "Bridgwater Bay",
    "Severn Estuary",
    "Axe Estuary",
    "Otter Estuary",
    "Exe Estuary",
    "Teign Estuary",
    "Dart Estuary",
    "Salcombe & Kingsbridge Estuary"
)

$hmhw = @(2.5945477590820385, 2.6007805710386314, 4.6442439823656851, 4.3084554248826956, 11.571311126048059, 0.62011173184357538, 1.7785714285714285, 2.5911602209944751, 2.6858638743455496, 4.8039215686274508, 1.4)
$hmtl = @(18.138109118982019, 5.0278214742332814, 4.5542310462498845, 3.4238383107056145, 9.481181257127977, 0.4351598173515982, 1.2589285714285714, 2.1226415094339623, 2.407766990291262, 5.1791044776119399, 2.4729241877256318)
$hmlw = @(16.069336471565538, 2.8271043496454937, 4.5509584694617597, 0.35342603388645066, 5.5095583550507792, 0.1673704414587332, 0.57796257796257799, 0.70107238605898126, 3.2123287671232879, 3.484320557491289, 1.3172413793103448)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws2.Range("A$row").Value = $names[$i]
    $ws2.Range("A$row").Font.Name = "Arial"

    $ws2.Range("B$row").Value = $hmhw[$i]
    $ws2.Range("C$row").Value = $hmtl[$i]
    $ws2.Range("D$row").Value = $hmlw[$i]
}

$ws2.Range("B2:D12").NumberFormat = "0.00"

# --- Sheet2 view tweaks ---------------------------------------------------
$ws2.Range("J23").Select()
